$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first sheet
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F25").Value = 2825
$wsExhibition.Range("F29").Value = 483
$wsExhibition.Range("F32").Value = 2293
$wsExhibition.Range("F45").Value = 43
$wsExhibition.Range("F46").Value = 504

# Sheet "全部类型" (All types) - fourth sheet, rows offset by +1
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F26").Value = 2825
$wsAll.Range("F30").Value = 483
$wsAll.Range("F33").Value = 2293
$wsAll.Range("F46").Value = 43
$wsAll.Range("F47").Value = 504
